$wb = $excel.ActiveWorkbook

# --- Service Contacts: move delivery_organisation_path column (old col R)
# --- to sit right after episode_key, before practitioner_key (new col D).
$ws = $wb.Worksheets.Item("Service Contacts")
$ws.Activate() | Out-Null

# Remember the custom width of the column being relocated so it travels
# along with its data (column width is a column-level property, not a
# per-cell one, so a plain cut/paste of cells won't carry it).
$savedWidth = $ws.Columns("R:R").ColumnWidth

# Insert a blank column before the practitioner_key column (D), shifting
# practitioner_key..service_contact_tags (old D:S) one column to the right
# (new E:T).
$ws.Columns("D:D").Insert() | Out-Null

# After the insert, the delivery_organisation_path column (formerly R) is
# now at S. Cut it and drop it into the newly inserted column D.
$ws.Columns("S:S").Cut($ws.Range("D1")) | Out-Null

# Remove the now-empty column left behind at S.
$ws.Columns("S:S").Delete() | Out-Null

# Re-apply the custom column width to its new home.
$ws.Columns("D:D").ColumnWidth = $savedWidth

# Match the recorded selection/scroll state for this sheet.
$ws.Columns("D:D").Select() | Out-Null

# --- Cosmetic selection-only changes on the other sheets ---

$wsOrg = $wb.Worksheets.Item("Organisations")
$wsOrg.Activate() | Out-Null
$wsOrg.Range("H1:J3").Select() | Out-Null

$wsK10 = $wb.Worksheets.Item("K10+")
$wsK10.Activate() | Out-Null
$wsK10.Range("F1:F5").Select() | Out-Null

$wsK5 = $wb.Worksheets.Item("K5")
$wsK5.Activate() | Out-Null
$wsK5.Range("F1:F5").Select() | Out-Null

# Restore Organisations as the active tab (it was active originally and
# the diff doesn't change that).
$wsOrg.Activate() | Out-Null
